$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.618.76'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.038.76'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.68%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '384.40'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.18%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.81'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.31%  '

$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.91'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.45%  '

$ws.Range('E11').Value = '  +0.10%  '

$ws.Range('E12').Value = '  +0.92%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.514.27'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.65%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.71'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.01%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.20%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.034.74'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.55%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.977'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.41%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.78'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -9.56%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '51.659.69'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.07%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.10'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.35%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.43'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.16%  '

$ws.Range('E22').Value = '  +0.13%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.17'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.30%  '

$ws.Range('E25').Value = '  -3.40%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.31'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +5.44%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.57'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.43%  '

$ws.Range('E28').Value = '  +4.34%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '26.38'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.00%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.108'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.08%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.31'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.02%  '

$ws.Range('E33').Value = '  +0.45%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '34.12'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.08%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.57'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.82%  '

$ws.Range('E36').Value = '  +3.08%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.35'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.47%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.290'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +7.34%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.10'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.13%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.87'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.65%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.117'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.12%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '127.74'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.58%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.54'
$ws.Range('D44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.69'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.90%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '21.71'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.18%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.47'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.54%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.08'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.92%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.034.25'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.88%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.338.87'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.79%  '

$ws.Range('E51').Value = '  +8.98%  '
